# ArrayExpress - Plant sample: add "Transcriptomics" tag to the TAGS block.
# Row 13: Tags | Plant | sample | Genomics | Transcriptomics(NEW) | mandatory(shifted from E->F)
# Row 14: Tags Term Accession Number | ...NCIT_C14258 | ...MS_1000457 | ...NCIT_C84343 | ...NCIT_C153189(NEW)
# Row 15: Tags Term Source REF | NCIT | MS | NCIT | NCIT(NEW)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("isa_template")

# Shift the existing "mandatory" tag from column E to column F on row 13.
$ws.Range("F13").Value = $ws.Range("E13").Value2

# Insert the new "Transcriptomics" tag in column E of row 13.
$ws.Range("E13").Value = "Transcriptomics"

# Add the accession number for the new tag on row 14.
$ws.Range("E14").Value = "http://purl.obolibrary.org/obo/NCIT_C153189"

# Add the term source REF for the new tag on row 15.
$ws.Range("E15").Value = "NCIT"
